# 6 Static Modules.pptx edit script
# 1) Refresh the two cached "today" date fields (Notes Master + slide layout 19)
#    from 7/9/2018 -> 7/16/2018.
# 2) Slide 5 ("Defining the Channels"): the picture that was showing the
#    (incorrect) "dynamic module" screenshot is swapped out - in practice this
#    means the existing picture placeholder is resized/repositioned and moved
#    to the front of the z-order (it now paints last, after every other shape
#    on the slide) and renamed to match the refreshed placeholder.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders -------------------------------------------------
# Notes Master "Date Placeholder 2" only accepts the change through the
# HeadersFooters facade in this host.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "7/16/2018"

# Slide layout 19's "Date Placeholder 3" (used by slides 2 & 3) only accepts
# the change when written straight through the shape's text range.
$slide2Layout = $p.Slides.Item(2).CustomLayout
$layoutDate = $slide2Layout.Shapes.Item(3)
$layoutDate.TextFrame.TextRange.Text = "7/16/2018"

# --- 2) Slide 5 picture ----------------------------------------------------
$s5 = $p.Slides.Item(5)
$pic = $s5.Shapes.Item(2)
$pic.Name = "Content Placeholder 4"

$EMU_PER_POINT = 914400 / 72
$pic.Left = 604838 / $EMU_PER_POINT
$pic.Top = 1806360 / $EMU_PER_POINT
$pic.Width = 4946850 / $EMU_PER_POINT
$pic.Height = 3886811 / $EMU_PER_POINT

$pic.ZOrder($msoBringToFront)
